# Delete the 2009年 data row (row 2). This shifts the 2010年 and 2011年
# rows up by one (to rows 2 and 3 respectively), matching the target
# layout where the sheet now only spans A1:P3 instead of A1:P4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()
